$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 220.63637
$ws.Range("I38").Value = 142.7
$ws.Range("J38").Value = 1000
$ws.Range("K38").Value = 428.1
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = -56.09999999999997
$ws.Range("N38").Value = -3744
# Row 39
$ws.Range("H39").Value = 479.86957
$ws.Range("I39").Value = 82.8
$ws.Range("J39").Value = 785.3077
$ws.Range("K39").Value = 248.4
$ws.Range("L39").Value = 2355.9231
$ws.Range("M39").Value = 47.60000000000002
$ws.Range("N39").Value = -2947.9231
# Row 62
$ws.Range("H62").Value = 2178.9375
$ws.Range("I62").Value = 1833.9166
$ws.Range("J62").Value = 3214
$ws.Range("K62").Value = 1833.9166
$ws.Range("L62").Value = 3214
$ws.Range("M62").Value = -1209.9166
$ws.Range("N62").Value = -4462
# Row 65
$ws.Range("H65").Value = 2178.9375
$ws.Range("I65").Value = 1833.9166
$ws.Range("J65").Value = 3214
$ws.Range("K65").Value = 9169.583000000001
$ws.Range("L65").Value = 16070
$ws.Range("M65").Value = -6049.583000000001
$ws.Range("N65").Value = -22310
# Row 103
$ws.Range("H103").Value = 1113
$ws.Range("I103").Value = 680.8
$ws.Range("J103").Value = 1833.3334
$ws.Range("K103").Value = 2042.4
$ws.Range("L103").Value = 5500.0002
$ws.Range("M103").Value = -1456.4
$ws.Range("N103").Value = -6672.0002
# Row 106
$ws.Range("H106").Value = 4005.8948
$ws.Range("I106").Value = 4306.6665
$ws.Range("J106").Value = 2878
$ws.Range("K106").Value = 4306.6665
$ws.Range("L106").Value = 2878
$ws.Range("M106").Value = -3675.6665
$ws.Range("N106").Value = -4140
# Row 127
$ws.Range("H127").Value = 1142.0322
$ws.Range("J127").Value = 1345.36
$ws.Range("L127").Value = 4036.08
$ws.Range("N127").Value = -13956.08
# Row 129
$ws.Range("H129").Value = 820.3200000000001
$ws.Range("J129").Value = 959.9737
$ws.Range("L129").Value = 2879.9211
$ws.Range("N129").Value = -12879.9211

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 4023.6
$ws.Range("I45").Value = 4030.5217
$ws.Range("J45").Value = 3944
$ws.Range("K45").Value = 4030.5217
$ws.Range("L45").Value = 3944
$ws.Range("M45").Value = -3653.5217
$ws.Range("N45").Value = -4698
# Row 63
$ws.Range("H63").Value = 2711.75
$ws.Range("I63").Value = 2711.75
$ws.Range("K63").Value = 2711.75
$ws.Range("M63").Value = -2025.75
# Row 64
$ws.Range("H64").Value = 22091
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 22091
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 22091
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -22587
# Row 66
$ws.Range("H66").Value = 2711.75
$ws.Range("I66").Value = 2711.75
$ws.Range("K66").Value = 13558.75
$ws.Range("M66").Value = -10126.75
# Row 67
$ws.Range("H67").Value = 22091
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 22091
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 22091
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -23807

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 62
$ws.Range("H62").Value = 32181
$ws.Range("J62").Value = 32181
$ws.Range("L62").Value = 32181
$ws.Range("N62").Value = -33553
# Row 65
$ws.Range("H65").Value = 32181
$ws.Range("J65").Value = 32181
$ws.Range("L65").Value = 96543
$ws.Range("N65").Value = -103407
# Row 99
$ws.Range("H99").Value = 1329
$ws.Range("I99").Value = 1223.75
$ws.Range("J99").Value = 1750
$ws.Range("K99").Value = 1223.75
$ws.Range("L99").Value = 1750
$ws.Range("M99").Value = 274.25
$ws.Range("N99").Value = -4746

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 326.33334
$ws.Range("I22").Value = 290.45456
$ws.Range("J22").Value = 425
$ws.Range("K22").Value = 290.45456
$ws.Range("L22").Value = 425
$ws.Range("M22").Value = 59.54543999999999
$ws.Range("N22").Value = -1125
# Row 31
$ws.Range("H31").Value = 10361.27
$ws.Range("I31").Value = 31928.223
$ws.Range("J31").Value = 1734.4889
$ws.Range("K31").Value = 31928.223
$ws.Range("L31").Value = 1734.4889
$ws.Range("M31").Value = -31633.223
$ws.Range("N31").Value = -2324.4889
# Row 34
$ws.Range("H34").Value = 10361.27
$ws.Range("I34").Value = 31928.223
$ws.Range("J34").Value = 1734.4889
$ws.Range("K34").Value = 31928.223
$ws.Range("L34").Value = 1734.4889
$ws.Range("M34").Value = -31726.223
$ws.Range("N34").Value = -2138.4889

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 5504557
$ws.Range("I11").Value = 7899472.5
$ws.Range("J11").Value = 714726
$ws.Range("K11").Value = 23698417.5
$ws.Range("L11").Value = 2144178
$ws.Range("M11").Value = -23698277.5
$ws.Range("N11").Value = -2144458
# Row 12
$ws.Range("H12").Value = 552.619
$ws.Range("I12").Value = 629.1177
$ws.Range("J12").Value = 227.5
$ws.Range("K12").Value = 1887.3531
$ws.Range("L12").Value = 682.5
$ws.Range("M12").Value = -1714.3531
$ws.Range("N12").Value = -1028.5
# Row 68
$ws.Range("H68").Value = 875.625
$ws.Range("I68").Value = 470.92856
$ws.Range("K68").Value = 1412.78568
$ws.Range("M68").Value = -601.78568
# Row 71
$ws.Range("H71").Value = 875.625
$ws.Range("I71").Value = 470.92856
$ws.Range("K71").Value = 4238.35704
$ws.Range("M71").Value = -182.3570399999999
# Row 100
$ws.Range("H100").Value = 3410
$ws.Range("J100").Value = 3410
$ws.Range("L100").Value = 10230
$ws.Range("N100").Value = -11852
# Row 131
$ws.Range("H131").Value = 795.3929000000001
$ws.Range("I131").Value = 503.625
$ws.Range("J131").Value = 912.1
$ws.Range("K131").Value = 1510.875
$ws.Range("L131").Value = 2736.3
$ws.Range("M131").Value = 3529.125
$ws.Range("N131").Value = -12816.3
# Row 132
$ws.Range("H132").Value = 2341.524
$ws.Range("I132").Value = 2079.111
$ws.Range("J132").Value = 2538.3333
$ws.Range("K132").Value = 18711.999
$ws.Range("L132").Value = 22844.9997
$ws.Range("M132").Value = -16181.999
$ws.Range("N132").Value = -27904.9997

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 45
$ws.Range("H45").Value = 26963
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
# Row 70
$ws.Range("H70").Value = 43519.23
$ws.Range("I70").Value = 67431.25
$ws.Range("J70").Value = 5260
$ws.Range("K70").Value = 67431.25
$ws.Range("L70").Value = 5260
$ws.Range("M70").Value = -67161.25
$ws.Range("N70").Value = -5800
# Row 73
$ws.Range("H73").Value = 43519.23
$ws.Range("I73").Value = 67431.25
$ws.Range("J73").Value = 5260
$ws.Range("K73").Value = 67431.25
$ws.Range("L73").Value = 5260
$ws.Range("M73").Value = -66495.25
$ws.Range("N73").Value = -7132
# Row 126
$ws.Range("H126").Value = 3774.7144
$ws.Range("I126").Value = 1755
$ws.Range("J126").Value = 4582.6
$ws.Range("K126").Value = 5265
$ws.Range("L126").Value = 13747.8
$ws.Range("M126").Value = -2795
$ws.Range("N126").Value = -18687.8
